# Update Name of Algo
# Apply updated KNN-imputed values in column D for the specified rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    11 = -7.084000000000001
    12 = -7.208
    15 = -8.297000000000001
    27 = -8.215999999999999
    28 = -8.108000000000001
    31 = -8.009
    32 = -7.753
    36 = -8.051
    38 = -7.662000000000001
    46 = -7.995
    54 = -8.247
    55 = -8.129000000000001
    56 = -8.17
    67 = -7.281000000000001
    69 = -7.321000000000001
    72 = -7.434
    73 = -8.029000000000002
    83 = -7.932
    86 = -7.885999999999998
    91 = -6.863000000000001
    93 = -7.568
    99 = -8.029999999999999
}

foreach ($row in $updates.Keys) {
    $ws.Range("D$row").Value = $updates[$row]
}
